# Applies the changes described by the diff:
#  - Highlights row 14 (A14:V14) with a solid yellow fill (adds a new fill + cellXf)
#  - Updates the importer name in P14 and P32
#  - Updates a batch of recalculated currency-conversion figures (T/U/V columns)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text value changes -----------------------------------------------
$ws.Range("P14").Value2 = "zuari cement limited"
$ws.Range("P32").Value2 = "continental india private limited"

# --- Highlight row 14 with solid yellow fill ---------------------------
$ws.Range("A14:V14").Interior.Color = 65535

# --- Recalculated numeric values ---------------------------------------
$ws.Range("T12").Value2 = 2.2005
$ws.Range("U12").Value2 = 2860.6286
$ws.Range("V12").Value2 = 3.4818

$ws.Range("T15").Value2 = 6.3817
$ws.Range("U15").Value2 = 6381471.8322

$ws.Range("T16").Value2 = 6.3817
$ws.Range("U16").Value2 = 12762943.3922

$ws.Range("T17").Value2 = 6.3817
$ws.Range("U17").Value2 = 15953678.8999

$ws.Range("T18").Value2 = 6.3817
$ws.Range("U18").Value2 = 1749448.5212

$ws.Range("T19").Value2 = 6.3817
$ws.Range("U19").Value2 = 3190735.78

$ws.Range("T31").Value2 = 2.1819
$ws.Range("U31").Value2 = 3491.0509
$ws.Range("V31").Value2 = 3.4818

$ws.Range("T34").Value2 = 84.8843
$ws.Range("U34").Value2 = 848818.9044999999
$ws.Range("V34").Value2 = 0.8065

$ws.Range("T35").Value2 = 82.98909999999999
$ws.Range("U35").Value2 = 331956.3307
$ws.Range("V35").Value2 = 0.8065

$ws.Range("T36").Value2 = 82.98909999999999
$ws.Range("U36").Value2 = 331956.3307
$ws.Range("V36").Value2 = 0.8065

$ws.Range("T37").Value2 = 671.5915
$ws.Range("U37").Value2 = 134317.2859
$ws.Range("V37").Value2 = 8.2247

$ws.Range("T38").Value2 = 79873.7458
$ws.Range("U38").Value2 = 1916969.9
$ws.Range("V38").Value2 = 984.2729

$ws.Range("T43").Value2 = 262171.2558
$ws.Range("U43").Value2 = 1310856.2557
$ws.Range("V43").Value2 = 2655.7896

Write-Host "edit complete"
